$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("client")

$ws.Range("E1").Value = "creation"
$ws.Range("E2").Value = "'2024-07-22 12:34:56"
$ws.Range("E2").NumberFormat = "m/d/yy h:mm"

$ws.Range("E3").Select()
